$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Localized quality-name strings were re-translated from Spanish to
# Portuguese: "Avanzado" -> "Avançado", "Legendario" -> "Lendário",
# "Desconocido" -> "Um estranho". Update every cell in column C that
# currently holds one of the old values.

$avancadoCells    = @("C2","C3","C19","C24","C65","C66","C67","C68","C69")
$legendarioCells  = @("C4","C20")
$desconocidoCells = @("C10","C14")

foreach ($c in $avancadoCells) {
    $ws.Range($c).Value = "Avançado"
}
foreach ($c in $legendarioCells) {
    $ws.Range($c).Value = "Lendário"
}
foreach ($c in $desconocidoCells) {
    $ws.Range($c).Value = "Um estranho"
}

# Refresh the sheet's saved selection / scroll position to match the
# author's final view: fully selecting the translated column (C2:C69)
# with no frozen/offset top-left cell.
$ws.Range("C2:C69").Select()
